# Updates cryptos list (price/volume figures) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to plain Text before writing so Excel does not
# reinterpret numeric-looking strings (e.g. "233.79") as numbers, matching
# the source data which stores every Coin/Link/Price/Volume cell as text.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "37.779.69"
Set-TextValue "E2" "  -0.17%  "
Set-TextValue "D3" "2.081.94"
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "233.79"
Set-TextValue "E5" "  -0.14%  "
Set-TextValue "E6" "  +0.08%  "
Set-TextValue "D7" "58.69"
Set-TextValue "E7" "  -0.74%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.03%  "
Set-TextValue "E9" "  +0.56%  "
Set-TextValue "D10" "0.0787"
Set-TextValue "E10" "  -0.25%  "
Set-TextValue "E11" "  +3.23%  "
Set-TextValue "D12" "15.05"
Set-TextValue "E12" "  +2.31%  "
Set-TextValue "D13" "2.388.90"
Set-TextValue "E13" "  -0.14%  "
Set-TextValue "D14" "21.48"
Set-TextValue "E14" "  +0.91%  "
Set-TextValue "E15" "  +1.84%  "
Set-TextValue "E16" "  +2.00%  "
Set-TextValue "D17" "2.081.44"
Set-TextValue "E17" "  -0.20%  "
Set-TextValue "D18" "37.778.87"
Set-TextValue "E18" "  +0.04%  "
Set-TextValue "D19" "6.14"
Set-TextValue "E19" "  -0.71%  "
Set-TextValue "D20" "71.45"
Set-TextValue "E20" "  +0.14%  "
Set-TextValue "E21" "  +1.36%  "
Set-TextValue "D22" "230.51"
Set-TextValue "E22" "  +0.70%  "
Set-TextValue "E23" "  -0.07%  "
Set-TextValue "E24" "  -0.73%  "
Set-TextValue "D25" "2.41"
Set-TextValue "E25" "  +1.44%  "
Set-TextValue "D26" "9.95"
Set-TextValue "E26" "  +10.60%  "
Set-TextValue "D27" "172.13"
Set-TextValue "E27" "  +1.11%  "
Set-TextValue "D28" "0.137"
Set-TextValue "E28" "  -1.21%  "
Set-TextValue "E29" "  -0.04%  "
Set-TextValue "E30" "  -0.13%  "
Set-TextValue "E31" "  +1.41%  "
Set-TextValue "D32" "4.75"
Set-TextValue "E32" "  +0.99%  "
Set-TextValue "E33" "  +1.06%  "
Set-TextValue "E34" "  -1.73%  "
Set-TextValue "E35" "  -1.17%  "
Set-TextValue "B36" "RenderToken"
Set-TextValue "C36" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D36" "3.41"
Set-TextValue "E36" "  -0.96%  "
Set-TextValue "B37" "WEMIXToken"
Set-TextValue "C37" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D37" "1.82"
Set-TextValue "E37" "  -0.77%  "
Set-TextValue "E38" "  +0.04%  "
Set-TextValue "E39" "  +0.79%  "
Set-TextValue "D40" "0.0235"
Set-TextValue "E40" "  +9.39%  "
Set-TextValue "D41" "102.53"
Set-TextValue "E41" "  +4.07%  "
Set-TextValue "D42" "0.0975"
Set-TextValue "E42" "  -1.67%  "
Set-TextValue "D43" "2.93"
Set-TextValue "E43" "  -1.15%  "
Set-TextValue "D44" "16.83"
Set-TextValue "E44" "  +4.55%  "
Set-TextValue "D45" "1.449.59"
Set-TextValue "E45" "  -0.65%  "
Set-TextValue "E46" "  -0.71%  "
Set-TextValue "D48" "4.09"
Set-TextValue "E48" "  -8.86%  "
Set-TextValue "E49" "  -0.79%  "
Set-TextValue "E50" "  -1.63%  "
Set-TextValue "E51" "  -0.11%  "
